# Inserts a new weekly price record for "Jengibre" (Vega Modelo de Temuco)
# as row 32, pushing all the existing records (previously rows 32-136) down
# by one row (to rows 33-137).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at position 32; Excel shifts rows 32..136 down to 33..137
# and carries the row-above formatting (e.g. the date style on column D).
$ws.Rows("32:32").Insert()

# Populate the newly inserted row 32 with the new observation.
$ws.Cells.Item(32, 1).Value = 10
$ws.Cells.Item(32, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(32, 3).Value = "La Araucanía"
$ws.Cells.Item(32, 4).Value = 44623
$ws.Cells.Item(32, 5).Value = 9
$ws.Cells.Item(32, 6).Value = 100114007
$ws.Cells.Item(32, 7).Value = "Jengibre"
$ws.Cells.Item(32, 8).Value = "Sin especificar"
$ws.Cells.Item(32, 9).Value = "Primera"
$ws.Cells.Item(32, 10).Value = 160
$ws.Cells.Item(32, 11).Value = 22000
$ws.Cells.Item(32, 12).Value = 25000
$ws.Cells.Item(32, 13).Value = 23500
$ws.Cells.Item(32, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(32, 15).Value = "Perú"
$ws.Cells.Item(32, 16).Value = 1808
$ws.Cells.Item(32, 17).Value = 13
$ws.Cells.Item(32, 18).Value = "Hortaliza"
